$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format temporarily so numeric-looking strings
# (e.g. "0.9939", "333.49") are not auto-converted to numbers by Excel.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.091.24"
$ws.Range("E2").Value = "  +3.52%  "
$ws.Range("D3").Value = "1.786.01"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("D5").Value = "335.43"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").Value = "0.9974"
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("D7").Value = "0.3821"
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").Value = "0.3434"
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("D9").Value = "47.69"
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("D10").Value = "1.151"
$ws.Range("E10").Value = "  -2.98%  "
$ws.Range("D11").Value = "0.07449"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").Value = "23.12"
$ws.Range("E12").Value = "  +7.08%  "
$ws.Range("D13").Value = "0.9987"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").Value = "6.405"
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D15").Value = "1.786.37"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").Value = "7.148"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").Value = "0.06648"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "82.90"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("D20").Value = "0.9938"
$ws.Range("E20").Value = "  -0.79%  "
$ws.Range("D21").Value = "17.52"
$ws.Range("E21").Value = "  +1.09%  "
$ws.Range("D22").Value = "6.451"
$ws.Range("E22").Value = "  -1.00%  "
$ws.Range("D23").Value = "28.094.89"
$ws.Range("E23").Value = "  +3.58%  "
$ws.Range("D24").Value = "12.08"
$ws.Range("E24").Value = "  -1.00%  "
$ws.Range("D25").Value = "2.389"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("D28").Value = "2.445"
$ws.Range("E28").Value = "  -1.97%  "
$ws.Range("D29").Value = "154.35"
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("D30").Value = "1.991.85"
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("D31").Value = "134.65"
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("D32").Value = "6.165"
$ws.Range("E32").Value = "  +2.40%  "
$ws.Range("D33").Value = "3.950"
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("D34").Value = "0.08734"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("D35").Value = "12.81"
$ws.Range("E35").Value = "  -1.87%  "
$ws.Range("D36").Value = "0.02428"
$ws.Range("E36").Value = "  +4.87%  "
$ws.Range("D37").Value = "0.6874"
$ws.Range("E37").Value = "  +0.68%  "
$ws.Range("D38").Value = "5.337"
$ws.Range("E38").Value = "  -0.81%  "
$ws.Range("D39").Value = "0.06343"
$ws.Range("E39").Value = "  +1.09%  "
$ws.Range("D40").Value = "0.2198"
$ws.Range("E40").Value = "  +1.33%  "
$ws.Range("D41").Value = "1.516"
$ws.Range("E41").Value = "  -6.68%  "
$ws.Range("D42").Value = "1.245"
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("D43").Value = "8.357"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("D44").Value = "14.34"
$ws.Range("E44").Value = "  +1.34%  "
$ws.Range("D45").Value = "0.9951"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("D46").Value = "0.6327"
$ws.Range("E46").Value = "  -1.81%  "
$ws.Range("D47").Value = "3.842"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").Value = "132.38"
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("D49").Value = "2.097"
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("D50").Value = "0.07470"
$ws.Range("E50").Value = "  +5.37%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "20.91"
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "1.442"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "1.273"
$ws.Range("E51").Value = "  +8.08%  "

# Restore column D cell style so no residual number-format metadata remains.
$dRange.Style = "Normal"

